$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5

# Row 3 updates
$ws.Range("Q3").Value = 1.65
$ws.Range("R3").Value = 2.2
